$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 32.363636
$ws.Range("I4").Value = 32.363636
$ws.Range("K4").Value = 32.363636
$ws.Range("M4").Value = 81.636364
$ws.Range("H7").Value = 8253
$ws.Range("I7").Value = 9000
$ws.Range("J7").Value = 8004
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 8004
$ws.Range("M7").Value = -8888
$ws.Range("N7").Value = -8228
$ws.Range("H14").Value = 8253
$ws.Range("I14").Value = 9000
$ws.Range("J14").Value = 8004
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 8004
$ws.Range("M14").Value = -8809
$ws.Range("N14").Value = -8386
$ws.Range("H39").Value = 214.9
$ws.Range("I39").Value = 224.88461
$ws.Range("J39").Value = 150
$ws.Range("K39").Value = 674.65383
$ws.Range("L39").Value = 450
$ws.Range("M39").Value = -378.65383
$ws.Range("N39").Value = -1042
$ws.Range("H69").Value = 14016.333
$ws.Range("J69").Value = 18749.5
$ws.Range("L69").Value = 56248.5
$ws.Range("N69").Value = -57996.5
$ws.Range("H72").Value = 14016.333
$ws.Range("J72").Value = 18749.5
$ws.Range("L72").Value = 168745.5
$ws.Range("N72").Value = -177481.5
$ws.Range("H101").Value = 825.6316
$ws.Range("I101").Value = 686.6875
$ws.Range("J101").Value = 1566.6666
$ws.Range("K101").Value = 2060.0625
$ws.Range("L101").Value = 4699.9998
$ws.Range("M101").Value = -438.0625
$ws.Range("N101").Value = -7943.9998
$ws.Range("H103").Value = 365.75
$ws.Range("J103").Value = 154.33333
$ws.Range("L103").Value = 462.99999
$ws.Range("N103").Value = -1634.99999
$ws.Range("H129").Value = 44180.25
$ws.Range("I129").Value = 44180.25
$ws.Range("K129").Value = 132540.75
$ws.Range("M129").Value = -127540.75
$ws.Range("H131").Value = 627402.9
$ws.Range("J131").Value = 5552.5
$ws.Range("L131").Value = 16657.5
$ws.Range("N131").Value = -26737.5
$ws.Range("H132").Value = 4884.1943
$ws.Range("I132").Value = 3809.0154
$ws.Range("K132").Value = 11427.0462
$ws.Range("M132").Value = -8897.046200000001
$ws.Range("H137").Value = 12306.275
$ws.Range("I137").Value = 3320.8333
$ws.Range("K137").Value = 9962.499899999999
$ws.Range("M137").Value = -7412.499899999999
$ws.Range("H138").Value = 2275.09
$ws.Range("I138").Value = 963.9259
$ws.Range("J138").Value = 2760.041
$ws.Range("K138").Value = 2891.7777
$ws.Range("L138").Value = 8280.123
$ws.Range("M138").Value = 2248.2223
$ws.Range("N138").Value = -18560.123

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7034.9395
$ws.Range("I32").Value = 3912.875
$ws.Range("J32").Value = 15360.444
$ws.Range("K32").Value = 3912.875
$ws.Range("L32").Value = 15360.444
$ws.Range("M32").Value = -3625.875
$ws.Range("N32").Value = -15934.444
$ws.Range("H132").Value = 3862334
$ws.Range("I132").Value = 5210.1333
$ws.Range("K132").Value = 15630.3999
$ws.Range("M132").Value = -13100.3999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3796.8948
$ws.Range("I94").Value = 2665
$ws.Range("K94").Value = 2665
$ws.Range("M94").Value = -2214
$ws.Range("H99").Value = 18266
$ws.Range("I99").Value = 932.3333
$ws.Range("K99").Value = 932.3333
$ws.Range("M99").Value = 565.6667
$ws.Range("H107").Value = 3248.1738
$ws.Range("I107").Value = 3914.0667
$ws.Range("K107").Value = 3914.0667
$ws.Range("M107").Value = -1994.0667
$ws.Range("H134").Value = 7234.6
$ws.Range("I134").Value = 3470.1162
$ws.Range("J134").Value = 20724
$ws.Range("K134").Value = 10410.3486
$ws.Range("L134").Value = 62172
$ws.Range("M134").Value = -7875.348599999999
$ws.Range("N134").Value = -67242

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12526.2
$ws.Range("I16").Value = 7405.1
$ws.Range("K16").Value = 7405.1
$ws.Range("M16").Value = -7118.1
$ws.Range("H31").Value = 17138.035
$ws.Range("I31").Value = 10752.333
$ws.Range("K31").Value = 10752.333
$ws.Range("M31").Value = -10457.333
$ws.Range("H34").Value = 17138.035
$ws.Range("I34").Value = 10752.333
$ws.Range("K34").Value = 10752.333
$ws.Range("M34").Value = -10550.333
$ws.Range("H86").Value = 33636.2
$ws.Range("I86").Value = 33636.2
$ws.Range("K86").Value = 33636.2
$ws.Range("M86").Value = -32513.2
$ws.Range("H89").Value = 33636.2
$ws.Range("I89").Value = 33636.2
$ws.Range("K89").Value = 168181
$ws.Range("M89").Value = -162565
$ws.Range("H113").Value = 12526.2
$ws.Range("I113").Value = 7405.1
$ws.Range("K113").Value = 7405.1
$ws.Range("M113").Value = -5235.1
$ws.Range("H122").Value = 1618.2222
$ws.Range("I122").Value = 1678
$ws.Range("J122").Value = 1498.6666
$ws.Range("K122").Value = 5034
$ws.Range("L122").Value = 4495.9998
$ws.Range("M122").Value = -2584
$ws.Range("N122").Value = -9395.9998
$ws.Range("H132").Value = 4216.3037
$ws.Range("I132").Value = 1459.8334
$ws.Range("J132").Value = 9177.950000000001
$ws.Range("K132").Value = 4379.5002
$ws.Range("L132").Value = 27533.85
$ws.Range("M132").Value = -1849.5002
$ws.Range("N132").Value = -32593.85
$ws.Range("H140").Value = 249500
$ws.Range("J140").Value = 249500
$ws.Range("L140").Value = 249500
$ws.Range("N140").Value = -259860

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 391.375
$ws.Range("I2").Value = 291.73685
$ws.Range("J2").Value = 770
$ws.Range("K2").Value = 1750.4211
$ws.Range("L2").Value = 4620
$ws.Range("M2").Value = -1637.4211
$ws.Range("N2").Value = -4846
$ws.Range("H26").Value = 809.38464
$ws.Range("I26").Value = 2025.8
$ws.Range("J26").Value = 49.125
$ws.Range("K26").Value = 6077.4
$ws.Range("L26").Value = 147.375
$ws.Range("M26").Value = -5789.4
$ws.Range("N26").Value = -723.375
$ws.Range("H46").Value = 402.9091
$ws.Range("J46").Value = 261
$ws.Range("L46").Value = 783
$ws.Range("N46").Value = -965
$ws.Range("H107").Value = 1081902.1
$ws.Range("I107").Value = 4083.9285
$ws.Range("J107").Value = 2087865.8
$ws.Range("K107").Value = 12251.7855
$ws.Range("L107").Value = 6263597.4
$ws.Range("M107").Value = -10331.7855
$ws.Range("N107").Value = -6267437.4
$ws.Range("H131").Value = 1494.84
$ws.Range("J131").Value = 1494.84
$ws.Range("L131").Value = 4484.52
$ws.Range("N131").Value = -14564.52
$ws.Range("H133").Value = 7777.7144
$ws.Range("I133").Value = 6898.8
$ws.Range("K133").Value = 20696.4
$ws.Range("M133").Value = -15636.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 811.3889
$ws.Range("I107").Value = 387.91666
$ws.Range("K107").Value = 387.91666
$ws.Range("M107").Value = 1532.08334
$ws.Range("H122").Value = 2264233.8
$ws.Range("I122").Value = 2829667
$ws.Range("K122").Value = 8489001
$ws.Range("M122").Value = -8486551
$ws.Range("H126").Value = 1718444.1
$ws.Range("I126").Value = 3600856.5
$ws.Range("K126").Value = 10802569.5
$ws.Range("M126").Value = -10800099.5
$ws.Range("H136").Value = 19806.572
$ws.Range("J136").Value = 19806.572
$ws.Range("L136").Value = 59419.716
$ws.Range("N136").Value = -64519.716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 27029790
$ws.Range("I22").Value = 2598.5312
$ws.Range("J22").Value = 200003820
$ws.Range("K22").Value = 2598.5312
$ws.Range("L22").Value = 200003820
$ws.Range("M22").Value = -2303.5312
$ws.Range("N22").Value = -200004410
$ws.Range("H27").Value = 27029790
$ws.Range("I27").Value = 2598.5312
$ws.Range("J27").Value = 200003820
$ws.Range("K27").Value = 2598.5312
$ws.Range("L27").Value = 200003820
$ws.Range("M27").Value = -2491.5312
$ws.Range("N27").Value = -200004034
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 5225
$ws.Range("I38").Value = 5800
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 5800
$ws.Range("L38").Value = 3500
$ws.Range("M38").Value = -5327
$ws.Range("N38").Value = -4446
$ws.Range("H92").Value = 74999.664
$ws.Range("J92").Value = 74999.664
$ws.Range("L92").Value = 74999.664
$ws.Range("N92").Value = -79991.664
$ws.Range("H132").Value = 4686.744
$ws.Range("I132").Value = 2417.262
$ws.Range("J132").Value = 100005
$ws.Range("K132").Value = 7251.786
$ws.Range("L132").Value = 300015
$ws.Range("M132").Value = -4721.786
$ws.Range("N132").Value = -305075
